# Applies the two text edits described by the commit:
#  1. Slide 20 ("Aggregation Framework"), shape "Text Placeholder 2" (id 7):
#     "02-04_MongoDB - Queries (2) - Aggregation Framework"
#       -> "01-04_MongoDB - Queries (2) - Aggregation Framework"
#  2. Slide 28 ("Free Courses on MongoDB University"), shape "Title 1" (id 2):
#     "Free Courses on MongoDB University"
#       -> split into "Free Courses on " + "MongoDB University " (two runs)

$p = $ppt.ActivePresentation

# --- Edit 1 -----------------------------------------------------------
$slide20 = $p.Slides.Item(20)
$shape20 = $slide20.Shapes.Item(2)
$tr20 = $shape20.TextFrame.TextRange

$oldText20 = "02-04_MongoDB - Queries (2) - Aggregation Framework"
$newText20 = "01-04_MongoDB - Queries (2) - Aggregation Framework"

$fullText20 = $tr20.Text
$pos20 = $fullText20.IndexOf($oldText20)
if ($pos20 -ge 0) {
    $run20 = $tr20.Characters($pos20 + 1, $oldText20.Length)
    $run20.Text = $newText20
}

# --- Edit 2 -----------------------------------------------------------
$slide28 = $p.Slides.Item(28)
$shape28 = $slide28.Shapes.Item(1)
$tr28 = $shape28.TextFrame.TextRange

$target28 = "MongoDB University"
$fullText28 = $tr28.Text
$pos28 = $fullText28.IndexOf($target28)
if ($pos28 -ge 0) {
    $run28 = $tr28.Characters($pos28 + 1, $target28.Length)
    $run28.Text = "MongoDB University "
}
